$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E text values are stored as text, not auto-converted to numbers
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.561.05"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "1.839.22"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "314.47"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4248"
$ws.Range("E7").Value = "  -3.74%  "
$ws.Range("D8").Value = "0.3632"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").Value = "45.61"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "0.07252"
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("D11").Value = "0.8891"
$ws.Range("E11").Value = "  -5.26%  "
$ws.Range("D12").Value = "20.54"
$ws.Range("E12").Value = "  -4.31%  "
$ws.Range("D13").Value = "1.903.80"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "5.368"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "6.561"
$ws.Range("E15").Value = "  -2.25%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "78.10"
$ws.Range("E18").Value = "  -4.90%  "
$ws.Range("D19").Value = "0.000008807"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "15.50"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("D22").Value = "27.542.00"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("D23").Value = "4.982"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").Value = "10.51"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").Value = "2.078.98"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").Value = "155.32"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "18.39"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "5.199"
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("D30").Value = "117.32"
$ws.Range("E30").Value = "  +3.14%  "
$ws.Range("D31").Value = "1.816"
$ws.Range("E31").Value = "  +5.20%  "
$ws.Range("D32").Value = "0.08883"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").Value = "0.7752"
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").Value = "4.550"
$ws.Range("E34").Value = "  -6.46%  "
$ws.Range("D35").Value = "2.959"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").Value = "1.097"
$ws.Range("E36").Value = "  -6.81%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "0.05398"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").Value = "0.01916"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("D41").Value = "2.762"
$ws.Range("E41").Value = "  -8.93%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.5052"
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "6.808"
$ws.Range("E43").Value = "  -4.39%  "
$ws.Range("D44").Value = "0.1646"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").Value = "8.178"
$ws.Range("E45").Value = "  -6.55%  "
$ws.Range("D46").Value = "0.06616"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "10.31"
$ws.Range("E47").Value = "  -2.83%  "
$ws.Range("D48").Value = "0.4684"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").Value = "104.96"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "1.627"
$ws.Range("E51").Value = "  -3.03%  "

# Restore default style (no explicit style index) while keeping text content
$ws.Range("D2:E51").Style = "Normal"
